# Edit the main text in paper rmd
# Update two cells in the ERGM term table on Sheet1:
#   B1: "Configuration (term in ERGM)" -> "Configuration (ERGM term)"
#   A8: "Centralization" -> "Centralization/popularity"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "Configuration (ERGM term)"
$ws.Range("A8").Value = "Centralization/popularity"

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("B13").Select()

$wb.Save()
